# The parser was updated to use TokenIteratorFieldRewriterSplit, which
# splits a single w:r run covering "{m" (open-brace + first letter of the
# field) into two separate runs: "{" and "m". The same kind of split is
# applied at the end of the field, where the run covering ")}" becomes two
# runs: ")" and "}".
#
# We locate each two-character token with Find, then toggle a character
# formatting property (Bold) off-and-on across just the first character of
# the match. In this Word engine, assigning a run-level formatting property
# to a sub-range of an existing run forces that run to be split into two
# runs at the sub-range boundary (one run per distinct set of properties);
# toggling the property back to its original value afterwards removes the
# property again without re-merging the now-separate runs, leaving a clean
# run split with the original rPr (here just w:lang) preserved on both
# halves.

$d = $word.ActiveDocument

function Split-RunAfterFirstChar($needle) {
    $found = $d.Content.Duplicate
    $found.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $firstChar = $d.Range($found.Start, $found.Start + 1)
    $firstChar.Bold = $true
    $firstChar.Bold = $false
}

# "{m" -> "{" + "m"
Split-RunAfterFirstChar("{m")

# ")}" -> ")" + "}"
Split-RunAfterFirstChar(")}")
